$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.045.95'
$ws.Range('E2').Value = '  -0.58%  '
$ws.Range('D3').Value = '1.830.28'
$ws.Range('D4').Value = '0.9990'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '240.93'
$ws.Range('E5').Value = '  -0.86%  '
$ws.Range('D6').Value = '0.6222'
$ws.Range('E6').Value = '  -6.26%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '0.07558'
$ws.Range('E8').Value = '  +1.74%  '
$ws.Range('D9').Value = '44.66'
$ws.Range('E9').Value = '  +6.64%  '
$ws.Range('E10').Value = '  -0.59%  '
$ws.Range('E11').Value = '  -0.81%  '
$ws.Range('D12').Value = '0.07624'
$ws.Range('E12').Value = '  -1.89%  '
$ws.Range('D13').Value = '1.830.69'
$ws.Range('E13').Value = '  -0.16%  '
$ws.Range('D14').Value = '4.954'
$ws.Range('E14').Value = '  -0.96%  '
$ws.Range('D15').Value = '0.6637'
$ws.Range('E15').Value = '  -0.66%  '
$ws.Range('D16').Value = '82.11'
$ws.Range('E16').Value = '  -1.08%  '
$ws.Range('D17').Value = '0.000009026'
$ws.Range('E17').Value = '  +7.62%  '
$ws.Range('D18').Value = '5.992'
$ws.Range('E18').Value = '  -2.21%  '
$ws.Range('D19').Value = '29.059.01'
$ws.Range('E19').Value = '  -0.48%  '
$ws.Range('D20').Value = '2.079.09'
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('D21').Value = '224.51'
$ws.Range('E21').Value = '  -1.75%  '
$ws.Range('E22').Value = '  -1.34%  '
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').Value = '7.181'
$ws.Range('E24').Value = '  +0.42%  '
$ws.Range('D25').Value = '1.001'
$ws.Range('D26').Value = '159.27'
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').Value = '8.402'
$ws.Range('E27').Value = '  -2.60%  '
$ws.Range('D28').Value = '0.1358'
$ws.Range('E28').Value = '  -3.16%  '
$ws.Range('D29').Value = '17.83'
$ws.Range('E29').Value = '  -0.96%  '
$ws.Range('D30').Value = '1.495'
$ws.Range('E30').Value = '  -1.55%  '
$ws.Range('E31').Value = '  +1.34%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '4.044'
$ws.Range('E32').Value = '  -1.82%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').Value = '4.020'
$ws.Range('E33').Value = '  -0.75%  '
$ws.Range('D34').Value = '0.05211'
$ws.Range('E34').Value = '  -1.34%  '
$ws.Range('D35').Value = '1.836'
$ws.Range('E35').Value = '  -1.59%  '
$ws.Range('D36').Value = '1.154'
$ws.Range('E36').Value = '  +1.13%  '
$ws.Range('D37').Value = '0.7312'
$ws.Range('E37').Value = '  -2.00%  '
$ws.Range('E38').Value = '  -0.25%  '
$ws.Range('D39').Value = '1.268.52'
$ws.Range('E39').Value = '  -3.58%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = '2.748'
$ws.Range('E40').Value = '  +0.33%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = '0.01783'
$ws.Range('E41').Value = '  -0.94%  '
$ws.Range('D42').Value = '6.350'
$ws.Range('E42').Value = '  +7.64%  '
$ws.Range('D43').Value = '0.8918'
$ws.Range('E43').Value = '  -4.22%  '
$ws.Range('E44').Value = '  +0.19%  '
$ws.Range('D45').Value = '101.57'
$ws.Range('E45').Value = '  -1.04%  '
$ws.Range('D46').Value = '1.976.92'
$ws.Range('E46').Value = '  -0.07%  '
$ws.Range('D47').Value = '0.5116'
$ws.Range('E47').Value = '  -0.61%  '
$ws.Range('D48').Value = '63.28'
$ws.Range('E48').Value = '  +0.23%  '
$ws.Range('E49').Value = '  -0.87%  '
$ws.Range('D50').Value = '0.3962'
$ws.Range('E50').Value = '  -1.38%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '8.893'
$ws.Range('E51').Value = '  +0.67%  '
